# Enviar.xlsx — replace the sample WhatsApp contact list with a fresh
# "teste" placeholder table (Pessoa / Número / Mensagem), per the commit:
# "Enviar mensagens no WhatsApp, pegando contatos da tabela no excel."
#
# Row 1 (headers) is untouched. Rows 2-4 get new contact data.
# Fill order matters for shared-string interning (matches how the values
# were actually typed: Mensagem column first, then Número, then Pessoa),
# so we touch column C, then B, then A.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Mensagem column - same value on every row
$ws.Range("C2").Value2 = "teste"
$ws.Range("C3").Value2 = "teste"
$ws.Range("C4").Value2 = "teste"

# Número column - placeholder phone numbers. Leading apostrophe keeps
# these as text (quote-prefixed), matching the existing text-formatted
# style already applied to these cells.
$ws.Range("B2").Value2 = "'0000000000000"
$ws.Range("B3").Value2 = "'0000000000000"
$ws.Range("B4").Value2 = "'0000000000000"

# Pessoa column - the new contacts
$ws.Range("A2").Value2 = "wesley"
$ws.Range("A3").Value2 = "silva"
$ws.Range("A4").Value2 = "almeida"

# Reset the view back to 100% zoom and move the selection.
$excel.ActiveWindow.Zoom = 100
$ws.Range("E4").Select()
